# Applies the "add displayDate to excel templates" edit:
#  - Renames "Sheet1" -> "Template" and "date-info" -> "Guide"
#  - Template (sheet1): inserts mods:originInfo / mods:displayDate
#    structure (splits the old combined originInfo+dateCreated cells
#    into separate tag cells, renames "Year" -> "Date Created", and adds
#    a brand-new "Display Date" field) while keeping the rest of the
#    MODS-fragment row intact; the trailing "example value" cell shifts
#    along with the column it documents.
#  - Guide (sheet2): replaces the old EDTF cheat-sheet with a two-column
#    Field Label / Notes and Best Practices reference table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename sheets
# ---------------------------------------------------------------------
$wsTemplate = $wb.Worksheets.Item(1)
$wsGuide = $wb.Worksheets.Item(2)
$wsTemplate.Name = "Template"
$wsGuide.Name = "Guide"

# ---------------------------------------------------------------------
# 2. Template sheet (sheet1) - rebuild header row 2
# ---------------------------------------------------------------------
# column -> [ value, isBoldLabel ]
$row2 = @{
    "A"  = @('<object pid="', $false)
    "B"  = @("PID", $true)
    "C"  = @('"><datastream type="md_descriptive" operation="update"><mods:mods xmlns:mods="http://www.loc.gov/mods/v3" xmlns:xlink="http://www.w3.org/1999/xlink" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance">', $false)
    "D"  = @("<mods:titleInfo><mods:title>", $false)
    "E"  = @("Title", $true)
    "F"  = @("</mods:title></mods:titleInfo>", $false)
    "G"  = @('<mods:identifier displayLabel="Digital Folder Number" type="local"> ', $false)
    "H"  = @("Digital Folder number", $true)
    "I"  = @("</mods:identifier>", $false)
    "J"  = @('<mods:identifier displayLabel="Accession Number" type="local">', $false)
    "K"  = @("Accession Number", $true)
    "L"  = @("</mods:identifier>", $false)
    "M"  = @("<mods:originInfo>", $false)
    "N"  = @('<mods:dateCreated encoding="edtf">', $false)
    "O"  = @("Date Created", $true)
    "P"  = @("</mods:dateCreated>", $false)
    "Q"  = @("<mods:displayDate>", $false)
    "R"  = @("Display Date", $true)
    "S"  = @("</mods:displayDate>", $false)
    "T"  = @("</mods:originInfo>", $false)
    "U"  = @('<mods:relatedItem type="original"><mods:identifier displayLabel="Source Media Identifier">', $false)
    "V"  = @("Source Media Identifier", $true)
    "W"  = @("</mods:identifier></mods:relatedItem>", $false)
    "X"  = @('<mods:identifier displayLabel="Collection Number" type="local">', $false)
    "Y"  = @("Collection Number", $true)
    "Z"  = @("</mods:identifier>", $false)
    "AA" = @("</mods:mods></datastream></object>", $false)
}

$colOrder = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

foreach ($col in $colOrder) {
    $addr = $col + "2"
    $pair = $row2[$col]
    $wsTemplate.Range($addr).Value = $pair[0]
}

# Apply "Text" number format across the whole header row in one shot
# (mirrors selecting the row/columns and setting Format Cells -> Text).
$wsTemplate.Range("A2:AA2").NumberFormat = "@"
$wsTemplate.Range("A1").NumberFormat = "@"

# Bold the field-label placeholder cells.
foreach ($col in $colOrder) {
    $pair = $row2[$col]
    if ($pair[1]) {
        $addr = $col + "2"
        $wsTemplate.Range($addr).Font.Bold = $true
    }
}

# B2 (PID) keeps its larger 12pt bold styling already present in the file.
$wsTemplate.Range("B2").Font.Bold = $true
$wsTemplate.Range("B2").Font.Size = 12

# Row 3: the blank "example answer" cell under Source Media Identifier
# moves from the old Q column to the new V column.
$wsTemplate.Range("Q3").ClearContents()
$wsTemplate.Range("V3").Value = ""
$wsTemplate.Range("V3").NumberFormat = "@"
$wsTemplate.Range("V3").WrapText = $true

# Selection / view bookkeeping to match the saved file.
[void]$wsTemplate.Range("B3").Select()

# ---------------------------------------------------------------------
# 3. Guide sheet (sheet2) - replace the EDTF cheat sheet with the new
#    Field Label / Notes and Best Practices reference table.
# ---------------------------------------------------------------------
$wsGuide.UsedRange.ClearContents()
$wsGuide.UsedRange.ClearFormats()

$wsGuide.Columns.Item(1).ColumnWidth = 27.85546875
$wsGuide.Columns.Item(2).ColumnWidth = 63.7109375

$guideRows = @(
    @("Field Label", "Notes and Best Practices"),
    @("PID", "DCR Object ID"),
    @("Title", ""),
    @("Digital Folder number", ""),
    @("Accession Number", ""),
    @("Date Created", "Date of creation of the resource, encoded according to EDTF. See https://adminliveunc.sharepoint.com/sites/DigitalCollectionsDocumentation/SitePages/Date-metadata-in-DCR.aspx for help and additional resources"),
    @("Display Date", "Optional free text date field that allows a date to be represented in human-readable form. It is recommended that mods:displayDate be entered to complement mods:dateCreated in order to provide a human-readable equivalent to the EDTF date. If you are unable to provide an EDTF date, providing only a mods:displayDate is preferable to no date."),
    @("Source Media Identifier", ""),
    @("Collection Number", '5-digit archival collection number. "70096"')
)

for ($i = 0; $i -lt $guideRows.Count; $i++) {
    $r = $i + 1
    $label = $guideRows[$i][0]
    $notes = $guideRows[$i][1]
    $wsGuide.Cells.Item($r, 1).Value = $label
    if ($notes -ne "") {
        $wsGuide.Cells.Item($r, 2).Value = $notes
    }
}

# Header row (Field Label / Notes and Best Practices) is bold.
$wsGuide.Range("A1").Font.Bold = $true
$wsGuide.Range("B1").Font.Bold = $true

[void]$wsGuide.Range("B9").Select()

Write-Output "done"
